$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 additions (existing row, bug fix + new column) ---
$ws.Range("A2").Value = 104

# --- Row 3 (new record) ---
$ws.Range("A3").Value = 105
$ws.Hyperlinks.Add($ws.Range("X3"), "https://goo.gl/maps/si4ws6UzUJcWKGr8A")
$ws.Range("X3").Style = "Hyperlink"
$ws.Range("C3").Value = "street"
$ws.Range("D3").Value = "South America"
$ws.Range("E3").Value = "Peru"
$ws.Range("F3").Value = "Lima"
$ws.Range("G3").Value = "Lima Metropolitan Area"
$ws.Range("H3").Value = "Province of Lima"
$ws.Range("I3").Value = "Villa María del Triunfo"
$ws.Range("K3").Value = 15828
$ws.Range("L3").Value = -12.1825599217897
$ws.Range("M3").Value = -76.948206122631603
$ws.Range("N3").Value = 2018
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 23
$ws.Hyperlinks.Add($ws.Range("W3"), "https://www.openstreetmap.org/way/241226313", "map=18/-12.18260/-76.94901")
$ws.Range("W3").Value = "https://www.openstreetmap.org/way/241226313#map=18/-12.18260/-76.94901"
$ws.Range("W3").Style = "Hyperlink"
$ws.Range("B3").Value = "Avenida Salvador Allende"
$ws.Range("V3").Value = 1

# --- Row 4 (new record) ---
$ws.Range("A4").Value = 106
$ws.Range("B4").Value = "Avenida Salvador Allende"
$ws.Range("C4").Value = "street"
$ws.Range("D4").Value = "South America"
$ws.Range("E4").Value = "Peru"
$ws.Range("F4").Value = "Lima"
$ws.Range("G4").Value = "Lima Metropolitan Area"
$ws.Range("H4").Value = "Province of Lima"
$ws.Range("I4").Value = "Villa María del Triunfo"

# --- Row 2 new column Q ---
$ws.Range("Q2").Value = "abacq date posted"
# --- Row 3 oldest_known_source ---
$ws.Range("Q3").Value = "openstreetmap"

# --- Row 4 remaining fields ---
$ws.Range("J4").Value = "Vallecito"
$ws.Range("K4").Value = 15811
$ws.Range("L4").Value = -12.1428039968166
$ws.Range("M4").Value = -76.953055786907996
$ws.Range("N4").Value = 2013
$ws.Range("O4").Value = 5
$ws.Hyperlinks.Add($ws.Range("X4"), "https://goo.gl/maps/9gS4FS77yjSxQ36x9")
$ws.Range("X4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("W4"), "https://www.openstreetmap.org/way/111941989")
$ws.Range("W4").Style = "Hyperlink"
$ws.Range("Q4").Value = "google maps"
$ws.Range("V4").Value = 1

# --- Row 2 google_maps_link (added last) ---
$ws.Hyperlinks.Add($ws.Range("X2"), "https://goo.gl/maps/tC5TJgUhPoRYyHkV9")
$ws.Range("X2").Style = "Hyperlink"

# --- Selection moved to A5 after data entry ---
$ws.Range("A5").Select()
